$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 10900
$ws.Range("I28").Value = 216.16667
$ws.Range("K28").Value = 216.16667
$ws.Range("M28").Value = 268.83333

# Row 132
$ws.Range("H132").Value = 6953.476
$ws.Range("I132").Value = 7751.643
$ws.Range("K132").Value = 23254.929
$ws.Range("M132").Value = -20724.929

# Row 138
$ws.Range("H138").Value = 2364.6045
$ws.Range("I138").Value = 2612.25
$ws.Range("J138").Value = 2311.7734
$ws.Range("K138").Value = 7836.75
$ws.Range("L138").Value = 6935.3202
$ws.Range("M138").Value = -2696.75
$ws.Range("N138").Value = -17215.3202

# Row 141
$ws.Range("H141").Value = 8742.714
$ws.Range("I141").Value = 4023.75
$ws.Range("J141").Value = 15034.667
$ws.Range("K141").Value = 12071.25
$ws.Range("L141").Value = 45104.001
$ws.Range("M141").Value = -6891.25
$ws.Range("N141").Value = -55464.001

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 409054.75
$ws.Range("I32").Value = 457924.9
$ws.Range("K32").Value = 457924.9
$ws.Range("M32").Value = -457637.9

# Row 74
$ws.Range("H74").Value = 2327.8948
$ws.Range("I74").Value = 2031.4706
$ws.Range("J74").Value = 2567.8572
$ws.Range("K74").Value = 2031.4706
$ws.Range("L74").Value = 2567.8572
$ws.Range("M74").Value = -1157.4706
$ws.Range("N74").Value = -4315.8572

# Row 77
$ws.Range("H77").Value = 2327.8948
$ws.Range("I77").Value = 2031.4706
$ws.Range("J77").Value = 2567.8572
$ws.Range("K77").Value = 10157.353
$ws.Range("L77").Value = 12839.286
$ws.Range("M77").Value = -5789.353000000001
$ws.Range("N77").Value = -21575.286

# Row 105
$ws.Range("H105").Value = 100000
$ws.Range("J105").Value = 100000
$ws.Range("L105").Value = 100000
$ws.Range("N105").Value = -106988

# Row 110
$ws.Range("H110").Value = 1283.25
$ws.Range("I110").Value = 1309
$ws.Range("K110").Value = 1309
$ws.Range("M110").Value = 736

# Row 132
$ws.Range("H132").Value = 3844.0952
$ws.Range("I132").Value = 3808.7827
$ws.Range("J132").Value = 3886.842
$ws.Range("K132").Value = 11426.3481
$ws.Range("L132").Value = 11660.526
$ws.Range("M132").Value = -8896.348100000001
$ws.Range("N132").Value = -16720.526

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 926.05
$ws.Range("I94").Value = 724.17645
$ws.Range("K94").Value = 724.17645
$ws.Range("M94").Value = -273.17645

# Row 141
$ws.Range("H141").Value = 24951.5
$ws.Range("I141").Value = 24709
$ws.Range("J141").Value = 25000
$ws.Range("K141").Value = 24709
$ws.Range("L141").Value = 25000
$ws.Range("M141").Value = -19529
$ws.Range("N141").Value = -35360

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6089.3125
$ws.Range("I31").Value = 1694.0358
$ws.Range("J31").Value = 12242.7
$ws.Range("K31").Value = 1694.0358
$ws.Range("L31").Value = 12242.7
$ws.Range("M31").Value = -1399.0358
$ws.Range("N31").Value = -12832.7

# Row 34
$ws.Range("H34").Value = 6089.3125
$ws.Range("I34").Value = 1694.0358
$ws.Range("J34").Value = 12242.7
$ws.Range("K34").Value = 1694.0358
$ws.Range("L34").Value = 12242.7
$ws.Range("M34").Value = -1492.0358
$ws.Range("N34").Value = -12646.7

# Row 58
$ws.Range("H58").Value = 1503.3478
$ws.Range("I58").Value = 1202
$ws.Range("J58").Value = 1779.5834
$ws.Range("K58").Value = 1202
$ws.Range("L58").Value = 1779.5834
$ws.Range("M58").Value = -999
$ws.Range("N58").Value = -2185.5834

# Row 86
$ws.Range("H86").Value = 2631.9355
$ws.Range("I86").Value = 2578.8333
$ws.Range("K86").Value = 2578.8333
$ws.Range("M86").Value = -1455.8333

# Row 89
$ws.Range("H89").Value = 2631.9355
$ws.Range("I89").Value = 2578.8333
$ws.Range("K89").Value = 12894.1665
$ws.Range("M89").Value = -7278.166499999999

# Row 136
$ws.Range("H136").Value = 1503.3478
$ws.Range("I136").Value = 1202
$ws.Range("J136").Value = 1779.5834
$ws.Range("K136").Value = 3606
$ws.Range("L136").Value = 5338.7502
$ws.Range("M136").Value = -1056
$ws.Range("N136").Value = -10438.7502

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 535
$ws.Range("I5").Value = 476.875
$ws.Range("K5").Value = 1430.625
$ws.Range("M5").Value = -1318.625

# Row 39
$ws.Range("H39").Value = 1581.4286
$ws.Range("J39").Value = 1581.4286
$ws.Range("L39").Value = 4744.2858
$ws.Range("N39").Value = -5332.2858

# Row 75
$ws.Range("H75").Value = 2104.2
$ws.Range("I75").Value = 1733.3334
$ws.Range("J75").Value = 2154.7727
$ws.Range("K75").Value = 5200.0002
$ws.Range("L75").Value = 6464.3181
$ws.Range("M75").Value = -4202.0002
$ws.Range("N75").Value = -8460.3181

# Row 78
$ws.Range("H78").Value = 2104.2
$ws.Range("I78").Value = 1733.3334
$ws.Range("J78").Value = 2154.7727
$ws.Range("K78").Value = 15600.0006
$ws.Range("L78").Value = 19392.9543
$ws.Range("M78").Value = -10608.0006
$ws.Range("N78").Value = -29376.9543

# Row 110
$ws.Range("H110").Value = 11604.193
$ws.Range("J110").Value = 12382.25
$ws.Range("L110").Value = 37146.75
$ws.Range("N110").Value = -45326.75

# Row 122
$ws.Range("H122").Value = 9106.083000000001
$ws.Range("I122").Value = 348.33334
$ws.Range("K122").Value = 3135.00006
$ws.Range("M122").Value = -685.0000600000003

# Row 131
$ws.Range("H131").Value = 989.8857400000001
$ws.Range("J131").Value = 1047.375
$ws.Range("L131").Value = 3142.125
$ws.Range("N131").Value = -13222.125

# Row 135
$ws.Range("H135").Value = 535
$ws.Range("I135").Value = 476.875
$ws.Range("K135").Value = 4291.875
$ws.Range("M135").Value = -1756.875

# Row 139
$ws.Range("H139").Value = 2675.963
$ws.Range("J139").Value = 2947.158
$ws.Range("L139").Value = 8841.474
$ws.Range("N139").Value = -19121.474

# Row 141
$ws.Range("H141").Value = 4658
$ws.Range("I141").Value = 1345.3846
$ws.Range("K141").Value = 4036.1538
$ws.Range("M141").Value = 1143.8462

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5759.4688
$ws.Range("I70").Value = 5756.696
$ws.Range("J70").Value = 5766.5557
$ws.Range("K70").Value = 5756.696
$ws.Range("L70").Value = 5766.5557
$ws.Range("M70").Value = -5486.696
$ws.Range("N70").Value = -6306.5557

# Row 73
$ws.Range("H73").Value = 5759.4688
$ws.Range("I73").Value = 5756.696
$ws.Range("J73").Value = 5766.5557
$ws.Range("K73").Value = 5756.696
$ws.Range("L73").Value = 5766.5557
$ws.Range("M73").Value = -4820.696
$ws.Range("N73").Value = -7638.5557

# Row 99
$ws.Range("H99").Value = 20933.166
$ws.Range("I99").Value = 12800
$ws.Range("J99").Value = 24999.75
$ws.Range("K99").Value = 12800
$ws.Range("L99").Value = 24999.75
$ws.Range("M99").Value = -10554
$ws.Range("N99").Value = -29491.75

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 4004.5652
$ws.Range("I132").Value = 3438.6924
$ws.Range("J132").Value = 4740.2
$ws.Range("K132").Value = 10316.0772
$ws.Range("L132").Value = 14220.6
$ws.Range("M132").Value = -7786.0772
$ws.Range("N132").Value = -19280.6

# Row 136
$ws.Range("H136").Value = 7577985
$ws.Range("I136").Value = 1950
$ws.Range("J136").Value = 16669227
$ws.Range("K136").Value = 5850
$ws.Range("L136").Value = 50007681
$ws.Range("M136").Value = -3300
$ws.Range("N136").Value = -50012781

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 4445.2666
$ws.Range("I81").Value = 4294
$ws.Range("J81").Value = 4618.143
$ws.Range("K81").Value = 8588
$ws.Range("L81").Value = 9236.286
$ws.Range("M81").Value = -7527
$ws.Range("N81").Value = -11358.286

# Row 84
$ws.Range("H84").Value = 4445.2666
$ws.Range("I84").Value = 4294
$ws.Range("J84").Value = 4618.143
$ws.Range("K84").Value = 42940
$ws.Range("L84").Value = 46181.43
$ws.Range("M84").Value = -37636
$ws.Range("N84").Value = -56789.43

# Row 113
$ws.Range("H113").Value = 1070.8
$ws.Range("I113").Value = 1070.8
$ws.Range("K113").Value = 3212.4
$ws.Range("M113").Value = -1042.4

# Row 116
$ws.Range("H116").Value = 56990
$ws.Range("J116").Value = 56990
$ws.Range("L116").Value = 56990
$ws.Range("N116").Value = -66168

# Row 126
$ws.Range("H126").Value = 1666.6666
$ws.Range("I126").Value = 1687.5
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 5062.5
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -2592.5
$ws.Range("N126").Value = -9440

# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# Row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# Row 135
$ws.Range("H135").Value = 111715
$ws.Range("J135").Value = 111715
$ws.Range("L135").Value = 111715
$ws.Range("N135").Value = -121855

# Row 136
$ws.Range("H136").Value = 4869.143
$ws.Range("I136").Value = 4651.636
$ws.Range("K136").Value = 13954.908
$ws.Range("M136").Value = -11404.908

# Row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
